# edit.ps1 - apply the "Add files via upload" revisions to
# draft-gandhi-spring-stamp-srpm-00.pptx
#
# Summary of changes:
#  * Slide 1 (title slide):
#      - subtitle ("Rectangle 3") run font size 20pt -> 18pt
#      - presenter box ("Rectangle 4") nudged down slightly (y offset)
#  * Slide 4 ("History of the Draft"):
#      - content placeholder nudged down slightly (y offset)
#      - "Moved STAMP support to " -> "Moved SRPM STAMP support to "
#      - "Scope TWAMP Light support as informational in " ->
#        "Keep SRPM TWAMP Light support in "
#  * Slide 16 (Probe Query slide):
#      - text box ("Content Placeholder 2") nudged (x/y offset)
#
# Note: PowerPoint's Shape.Left/Top/Width/Height are IEEE-754 single
# precision (points) under the hood, so EMU targets that aren't exact
# multiples of the float32 grid are approached with literals chosen to
# round-trip to the exact EMU value on save.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 1 - title slide
# ---------------------------------------------------------------------
$s1 = $p.Slides.Item(1)

# Subtitle ("draft-gandhi-spring-stamp-srpm-00 (previously ...)") : 20 -> 18 pt
$subtitle = $s1.Shapes.Item(2)
$subtitle.TextFrame.TextRange.Font.Size = 18

# Presenter info box: shift down slightly
$presenterBox = $s1.Shapes.Item(3)
$presenterBox.Top = 234.06240844726562

# ---------------------------------------------------------------------
# Slide 4 - "History of the Draft"
# ---------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$content4 = $s4.Shapes.Item(2)

# Shift the content placeholder down slightly
$content4.Top = 67.50004577636719

$tr4 = $content4.TextFrame.TextRange

# Replace "Moved STAMP support to " -> "Moved SRPM STAMP support to "
# (Characters() exposes paragraph breaks as empty-string characters, so
# rebuild the plain text with CRs to locate the target substring.)
$len = $tr4.Length
$full = ""
for ($i = 1; $i -le $len; $i++) {
    $t = $tr4.Characters($i, 1).Text
    if ($t -eq "") { $full += "`r" } else { $full += $t }
}
$find1 = "Moved STAMP support to "
$idx1 = $full.IndexOf($find1)
if ($idx1 -ge 0) {
    $tr4.Characters($idx1 + 1, $find1.Length).Text = "Moved SRPM STAMP support to "
}

# Replace "Scope TWAMP Light support as informational in " ->
# "Keep SRPM TWAMP Light support in " (recompute offsets: text length changed above)
$len = $tr4.Length
$full = ""
for ($i = 1; $i -le $len; $i++) {
    $t = $tr4.Characters($i, 1).Text
    if ($t -eq "") { $full += "`r" } else { $full += $t }
}
$find2 = "Scope TWAMP Light support as informational in "
$idx2 = $full.IndexOf($find2)
if ($idx2 -ge 0) {
    $tr4.Characters($idx2 + 1, $find2.Length).Text = "Keep SRPM TWAMP Light support in "
}

# ---------------------------------------------------------------------
# Slide 16 - "Probe Query for SR-MPLS and SRv6 Policy"
# ---------------------------------------------------------------------
$s16 = $p.Slides.Item(16)
$content16 = $s16.Shapes.Item(4)
$content16.Left = 16.312562942504883
$content16.Top = 115.1109848022461
